# Messages order.xlsx - row edit: delete the "PlayersOrderMessage" row
# (row 8) from Sheet1; the author apparently pulled the player-order
# broadcast out of PlayerActor (commit: "Further refactoring in
# PlayerActor. Created RandomPlayerActor."), so the corresponding message
# no longer appears in the message table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row 8 (PlayersOrderMessage | Game | Player) and shift
# everything below it up by one row.
$ws.Rows(8).Delete()

# The sheet always spans down to row 14 (A1:C14); after the delete the
# used range shrinks to row 13, so touch row 14 to keep it materialized
# as a (blank) row, matching the original sheet shape.
$ws.Rows(14).RowHeight = 23.25

# Leave the selection where it ended up after the delete/scroll.
$ws.Range("B10").Select()
